$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.089056333333333
$ws.Range("H2").Value = 6.267169
$ws.Range("I2").Value = 0.7196603919224289
$ws.Range("J2").Value = 0.719660391922429
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.332796666666667
$ws.Range("N2").Value = 9.998390000000001
$ws.Range("Q2").Value = 6.962399984212222
$ws.Range("R2").Value = 62.66159985791001
$ws.Range("S2").Value = 0.7196603919224289
$ws.Range("T2").Value = 0.719660391922429

# Row 3 updates
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.8137799999999999
$ws.Range("H3").Value = 2.44134
$ws.Range("I3").Value = 0.280339608077571
$ws.Range("J3").Value = 0.280339608077571
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.332796666666667
$ws.Range("N3").Value = 9.998390000000001
$ws.Range("Q3").Value = 2.7121632714
$ws.Range("R3").Value = 24.4094694426
$ws.Range("S3").Value = 0.280339608077571
$ws.Range("T3").Value = 0.280339608077571
